$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add "done" marker for Fosso_Ghiaia row (row 10)
$ws.Range("E10").Value = "done"

# Update selection to reflect where the user last clicked (B7)
$ws.Range("B7").Select()
